# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / Correspond Handoff/Handback
# datetimes for the f6c6938f file's de-de and zh-cn handback rows, plus the
# de-de column on the Overview summary sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: de-de column (G) for f6c6938f row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 16:52:40"

# --- zh-cn sheet: f6c6938f row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-23 16:52:35"
$wsZhCn.Range("K3").Value = "2016-08-23 16:52:51"

# --- de-de sheet: f6c6938f row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-23 16:52:40"
$wsDeDe.Range("K3").Value = "2016-08-23 16:52:58"
